# Update capital structure database values for rows 2 and 3
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "G"  = -0.2262931034482759
    "H"  = -0.2262931034482759
    "I"  = -0.6573275862068966
    "J"  = -0.6573275862068966
    "K"  = -1.39
    "L"  = -0.2995689655172414
    "M"  = 0.718
    "N"  = 0.1237931034482759
    "O"  = -0.516546762589928
    "S"  = 0.718
    "T"  = 1
    "U"  = 1.13
    "V"  = 0.1948275862068966
    "W"  = 2.260162601626016
    "X"  = 0.0601501458654287
    "Y"  = 2.200012455760588
    "Z"  = 8.498168498168498
    "AA" = -5.586080586080587
    "AB" = 0.05846780085996194
    "AC" = -5.644548386940548
    "AD" = 0.38
    "AF" = 0.38
    "AG" = -0.7499999999999999
    "AH" = 0.06148867313915858
    "AI" = 0.6333333333333334
    "AJ" = -0.1485148514851485
    "AK" = 1.415094339622641
    "AN" = -0.1310344827586207
    "AP" = 0.2586206896551724
}

foreach ($row in @(2, 3)) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
    # Columns removed from the dataset (no longer have data in these rows)
    $ws.Range("AO$row").ClearContents()
    $ws.Range("AQ$row").ClearContents()
}
